$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.908.16'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '1.550.39'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '207.18'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '21.62'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0587'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0860'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Value = '1.771.82'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '1.541.78'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.72'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '61.79'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.911.42'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '215.02'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('D19').Value = '0.0₃0688'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.49'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.66'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.85'
$ws.Range('D27').Style = "Normal"
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0463'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('D33').Value = '1.406.08'
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('E35').Value = '  +3.51%  '
$ws.Range('E36').Value = '  +2.19%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.521'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.808'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '63.56'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').Value = '1.686.21'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.15'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0953'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.67%  '
